$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the two runs that read "...usin" / "g discrete..." into a single
#    run with the combined text "...using discrete...", and drop the
#    <w:bookmarkStart/bookmarkEnd name="_GoBack"/> pair that currently sits
#    between them.
# ---------------------------------------------------------------------------

# Locate the two runs' text precisely via Find so we don't depend on fixed
# character offsets.
$run1 = $d.Content
$run1.Find.Execute(" Well-versed in various algorithm design paradigms usin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$run1Start = $run1.Start
$run1End = $run1.End

$run2 = $d.Content
$run2.Find.Execute("g discrete mathematics. Experienced in Object", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$run2Start = $run2.Start
$run2End = $run2.End

$mergedText = " Well-versed in various algorithm design paradigms using discrete mathematics. Experienced in Object"

# Drop helper bookmarks around the pair of runs first, so the forced text
# edit below (needed to make the engine actually re-merge the two runs)
# cannot bleed out into neighbouring runs that share identical formatting.
$before = $d.Range($run1Start, $run1Start)
$d.Bookmarks.Add("_mergeBarrierBefore", $before) | Out-Null
$after = $d.Range($run2End, $run2End)
$d.Bookmarks.Add("_mergeBarrierAfter", $after) | Out-Null

# Remove the original _GoBack bookmark that sits between the two runs.
$d.Bookmarks.Item("_GoBack").Delete()

# Force an actual text mutation spanning the boundary between the two runs
# so the engine recombines them into a single run.
$joint = $d.Range($run1End, $run1End)
$joint.InsertBefore("Z")
$jointChar = $d.Range($run1End, $run1End + 1)
$jointChar.Text = ""

# Sanity: the merged span should now read exactly the combined text.
$check = $d.Range($run1Start, $run2End)
if ($check.Text -ne $mergedText) {
    $check.Text = $mergedText
}

# Remove the helper barrier bookmarks now that the merge has happened.
$d.Bookmarks.Item("_mergeBarrierBefore").Delete()
$d.Bookmarks.Item("_mergeBarrierAfter").Delete()

# ---------------------------------------------------------------------------
# 2) Re-add the _GoBack bookmark at the very start of the document (this is
#    where Word leaves it after the last edit position is "the top" of the
#    document).
# ---------------------------------------------------------------------------

$docStart = $d.Range(0, 0)
$docStart.InsertBefore("X")
$placeholder = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $placeholder) | Out-Null
$placeholderAgain = $d.Range(0, 1)
$placeholderAgain.Text = ""

Write-Output "edit complete"
